# Auto-generated edit script applying scheduled market-data refresh
# to the Odin_Profits leve-profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# --- ALC sheet, hunk 0 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H115").Value = 714.5
$ws.Range("I115").Value = 530.8570999999999
$ws.Range("J115").Value = 2000
$ws.Range("K115").Value = 1592.5713
$ws.Range("L115").Value = 6000
$ws.Range("M115").Value = -25.57129999999984
$ws.Range("N115").Value = -9134

# --- ALC sheet, hunk 1 ---
$ws.Range("H116").Value = 15877988
$ws.Range("I116").Value = 37038972
$ws.Range("J116").Value = 7251.25
$ws.Range("K116").Value = 37038972
$ws.Range("L116").Value = 7251.25
$ws.Range("M116").Value = -37035530
$ws.Range("N116").Value = -14135.25

# --- ALC sheet, hunk 2 ---
$ws.Range("H131").Value = 1806.7
$ws.Range("I131").Value = 1729.6666
$ws.Range("K131").Value = 5188.9998
$ws.Range("M131").Value = -148.9997999999996

# --- ALC sheet, hunk 3 ---
$ws.Range("H132").Value = 530786.75
$ws.Range("I132").Value = 653940.8
$ws.Range("J132").Value = 13539.8
$ws.Range("K132").Value = 1961822.4
$ws.Range("L132").Value = 40619.39999999999
$ws.Range("M132").Value = -1959292.4
$ws.Range("N132").Value = -45679.39999999999

# --- ARM sheet, hunk 4 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1427951
$ws.Range("I32").Value = 2906.3242
$ws.Range("K32").Value = 2906.3242
$ws.Range("M32").Value = -2619.3242

# --- ARM sheet, hunk 5 ---
$ws.Range("H45").Value = 1706.9375
$ws.Range("I45").Value = 1450.7858
$ws.Range("J45").Value = 3500
$ws.Range("K45").Value = 1450.7858
$ws.Range("L45").Value = 3500
$ws.Range("M45").Value = -1073.7858
$ws.Range("N45").Value = -4254

# --- ARM sheet, hunk 6 ---
$ws.Range("H61").Value = 3881.8594
$ws.Range("I61").Value = 2420.2632
$ws.Range("K61").Value = 2420.2632
$ws.Range("M61").Value = -2208.2632

# --- ARM sheet, hunk 7 ---
$ws.Range("H132").Value = 1412292.6
$ws.Range("I132").Value = 1687776.4
$ws.Range("K132").Value = 5063329.199999999
$ws.Range("M132").Value = -5060799.199999999

# --- ARM sheet, hunk 8 ---
$ws.Range("H136").Value = 3881.8594
$ws.Range("I136").Value = 2420.2632
$ws.Range("K136").Value = 7260.7896
$ws.Range("M136").Value = -4710.7896

# --- BSM sheet, hunk 9 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 4294.5
$ws.Range("I54").Value = 4294.5
$ws.Range("K54").Value = 4294.5
$ws.Range("M54").Value = -3810.5

# --- BSM sheet, hunk 10 ---
$ws.Range("H94").Value = 4793.12
$ws.Range("J94").Value = 10325
$ws.Range("L94").Value = 10325
$ws.Range("N94").Value = -11227

# --- BSM sheet, hunk 11 ---
$ws.Range("H132").Value = 98946.5
$ws.Range("J132").Value = 98946.5
$ws.Range("L132").Value = 98946.5
$ws.Range("N132").Value = -109066.5

# --- BSM sheet, hunk 12 ---
$ws.Range("H134").Value = 11011.143
$ws.Range("I134").Value = 10609.25
$ws.Range("J134").Value = 12297.2
$ws.Range("K134").Value = 31827.75
$ws.Range("L134").Value = 36891.60000000001
$ws.Range("M134").Value = -29292.75
$ws.Range("N134").Value = -41961.60000000001

# --- CRP sheet, hunk 13 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 812797.7
$ws.Range("I22").Value = 1488436.6
$ws.Range("J22").Value = 2030.9
$ws.Range("K22").Value = 1488436.6
$ws.Range("L22").Value = 2030.9
$ws.Range("M22").Value = -1488086.6
$ws.Range("N22").Value = -2730.9

# --- CRP sheet, hunk 14 ---
$ws.Range("H31").Value = 10644573
$ws.Range("I31").Value = 45473000
$ws.Range("J31").Value = 2553.5833
$ws.Range("K31").Value = 45473000
$ws.Range("L31").Value = 2553.5833
$ws.Range("M31").Value = -45472705
$ws.Range("N31").Value = -3143.5833

# --- CRP sheet, hunk 15 ---
$ws.Range("H34").Value = 10644573
$ws.Range("I34").Value = 45473000
$ws.Range("J34").Value = 2553.5833
$ws.Range("K34").Value = 45473000
$ws.Range("L34").Value = 2553.5833
$ws.Range("M34").Value = -45472798
$ws.Range("N34").Value = -2957.5833

# --- CRP sheet, hunk 16 ---
$ws.Range("H99").Value = 6947562
$ws.Range("I99").Value = 13892035
$ws.Range("J99").Value = 3088.875
$ws.Range("K99").Value = 13892035
$ws.Range("L99").Value = 3088.875
$ws.Range("M99").Value = -13890537
$ws.Range("N99").Value = -6084.875

# --- CRP sheet, hunk 17 ---
$ws.Range("H126").Value = 6947562
$ws.Range("I126").Value = 13892035
$ws.Range("J126").Value = 3088.875
$ws.Range("K126").Value = 41676105
$ws.Range("L126").Value = 9266.625
$ws.Range("M126").Value = -41673635
$ws.Range("N126").Value = -14206.625

# --- CRP sheet, hunk 18 ---
$ws.Range("H141").Value = 276768.16
$ws.Range("J141").Value = 302121.8
$ws.Range("L141").Value = 302121.8
$ws.Range("N141").Value = -312481.8

# --- CUL sheet, hunk 19 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 668.875
$ws.Range("J33").Value = 930.8
$ws.Range("L33").Value = 5584.799999999999
$ws.Range("N33").Value = -6150.799999999999

# --- CUL sheet, hunk 20 ---
$ws.Range("H51").Value = 2999
$ws.Range("J51").Value = 2999
$ws.Range("L51").Value = 8997
$ws.Range("N51").Value = -9917

# --- CUL sheet, hunk 21 ---
$ws.Range("H132").Value = 3507.7
$ws.Range("J132").Value = 4281.143
$ws.Range("L132").Value = 38530.287
$ws.Range("N132").Value = -43590.287

# --- CUL sheet, hunk 22 ---
$ws.Range("H137").Value = 2217.6843
$ws.Range("I137").Value = 1368.1818
$ws.Range("J137").Value = 3385.75
$ws.Range("K137").Value = 4104.5454
$ws.Range("L137").Value = 10157.25
$ws.Range("M137").Value = 995.4546
$ws.Range("N137").Value = -20357.25

# --- GSM sheet, hunk 23 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 9199.666999999999
$ws.Range("I113").Value = 3399.5
$ws.Range("J113").Value = 14999.833
$ws.Range("K113").Value = 3399.5
$ws.Range("L113").Value = 14999.833
$ws.Range("M113").Value = -1229.5
$ws.Range("N113").Value = -19339.833

# --- GSM sheet, hunk 24 ---
$ws.Range("H126").Value = 45467372
$ws.Range("I126").Value = 62505140
$ws.Range("J126").Value = 33333
$ws.Range("K126").Value = 187515420
$ws.Range("L126").Value = 99999
$ws.Range("M126").Value = -187512950
$ws.Range("N126").Value = -104939

# --- GSM sheet, hunk 25 ---
$ws.Range("H132").Value = 3905.4092
$ws.Range("I132").Value = 3810.5312
$ws.Range("J132").Value = 4158.4165
$ws.Range("K132").Value = 11431.5936
$ws.Range("L132").Value = 12475.2495
$ws.Range("M132").Value = -8901.5936
$ws.Range("N132").Value = -17535.2495

# --- LTW sheet, hunk 26 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

# --- WVR sheet, hunk 27 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4691.647
$ws.Range("I126").Value = 2145.2222
$ws.Range("J126").Value = 7556.375
$ws.Range("K126").Value = 6435.6666
$ws.Range("L126").Value = 22669.125
$ws.Range("M126").Value = -3965.6666
$ws.Range("N126").Value = -27609.125

# --- WVR sheet, hunk 28 ---
$ws.Range("H132").Value = 25408770
$ws.Range("I132").Value = 22227524
$ws.Range("K132").Value = 66682572
$ws.Range("M132").Value = -66680042
